$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: F1 becomes "Kulcsszavak (magyar)" and add G1 "Kulcsszavak (angol)"
$ws.Range("F1").Value = "Kulcsszavak (magyar)"

# Give the new G1 header the same (bold/centered) formatting as the other
# header cells before writing its text
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G1").Value = "Kulcsszavak (angol)"

# English keyword values extracted from the PDFs, one per row (row 1 is header)
$englishKeywords = @{
    2  = "microcontroller, WiFi, gyroscope, accelerometer"
    3  = "IoT, Home -automation, Raspberry Pi,  MQTT"
    4  = "temperature measurement, bluetooth, application, robot control."
    6  = "backgammon, image processing, OpenCV ,  Python"
    8  = "network, protocol, routing"
    9  = "Markov clustering, graph, sparse matrix, efficient algorithm, bioinformatics, protein sequence  Keywords: Markov clustering, graph, sparse matrix, efficient algorithm, bioinformatics, protein sequence   _____________________________________________________________________________    14"
    10 = "platform -independent, mobile application, Flutter, Spring Boot, event."
    11 = "digital signal processing, sound effects, ARM programming, Python simulation"
    12 = "deep learning, cardiovascular diseases, spectrogram, signal filtering, disease diagnosis.            Tartalom"
    13 = "mobile robots, voice control"
    14 = "decentralized exchange, smart contract, crypto, blockchain"
    16 = "Internship, Recommendation System , Cosine Similarity , Angular, Spring Boot"
    17 = "facial recognition, identiﬁcation, biometrics, access control system, secu-rity"
    19 = "chest X-ray, classification, convolutional neural network, transfer learn-ing, GoogleNet architecture"
    20 = "user interface , smart system , microprocessor"
    21 = "Baxter, ROS, ArUco, Opencv, Python, image processing, industrial and educational robotics"
    22 = "Mobile application, PowerPoint, laser pointer."
    23 = "IoT, Smart Energy meter, MQTT, Energy Management"
    24 = "industrial robot, remote control, gestures, Kinect       8"
    25 = "traffic penalty management, android app,webinterface8"
}

foreach ($row in $englishKeywords.Keys) {
    $ws.Cells.Item($row, 7).Value = $englishKeywords[$row]
}

# Rows where the English keyword cell exists but is empty (mirrors the blank
# Hungarian keyword cells already present for these rows)
$emptyRows = @(5, 7, 15, 18)
foreach ($row in $emptyRows) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.Font.Bold = $cell.Font.Bold
}
